$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting from the existing header cell (H1) to the new header cells
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# Data values for columns I (I0) and J (IF), rows 2-16
$dataI = @(8, 9, 4, 1, 1, 6, 5, 1, 1, 1, 1, 2, 4, 7, 6)
$dataJ = @(9, 9, 8, 6, 2, 6, 6, 4, 6, 6, 5, 6, 4, 8, 6)

for ($i = 0; $i -lt 15; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
